$d = $word.ActiveDocument

# --- First paragraph (the **ID__...__ID** placeholder line) ---
$p1 = $d.Paragraphs(1)

# Add a paragraph border (no visible line, just 5pt space on every edge)
# -> <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

# Left indent: 120 twips (6pt) -> 225 twips (11.25pt)
$p1.Format.LeftIndent = 11.25

# Update the placeholder id text in place (keeps the existing run's formatting)
$idRange = $p1.Range.Duplicate
$idRange.Find.Execute("**ID__AFFARS_5309_topic_2__ID**", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0) | Out-Null
$idRange.Text = "**ID__AFFARS_SUBPART_5309_1__ID**"

# The paragraph used to end with a separate run containing a single trailing
# space; drop it so the paragraph ends right after the id text.
$p1 = $d.Paragraphs(1)
$trailingSpace = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}
